$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 3
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Value = "FR"
$ws.Range("C5").Value = "????erere ee"

$ws.Range("C6").Select() | Out-Null
